$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.977.09"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.632.50"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0881"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.863.63"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.624.66"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "27.977.40"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -3.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "1.407.34"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  +8.19%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").Value = "1.774.66"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -0.32%  "
